$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.570.40"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.630.40"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("E6").Value = "  +2.09%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").Value = "1.856.12"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "1.628.96"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.52%  "

$ws.Range("D17").Value = "26.701.61"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.03%  "

$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.31%  "

$ws.Range("E29").Value = "  +1.19%  "

$ws.Range("E30").Value = "  -2.59%  "

$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("E32").Value = "  +2.85%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").Value = "1.217.75"
$ws.Range("E35").Value = "  +4.69%  "

$ws.Range("E36").Value = "  -1.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0173"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.44%  "

$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("E39").Value = "  -1.46%  "

$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("E41").Value = "  -2.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.796"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "

$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("D44").Value = "1.767.81"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
